$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country Updates")

# Column B holds the "date of report" for each row; the daily refresh
# bumps every occurrence of the previous day's serial date (43926,
# 2020-04-05) forward to the new day's serial date (43927, 2020-04-06).
$range = $ws.Range("B5:B74")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 43926) {
        $cell.Value = 43927
    }
}
